# Applies the "Updated: st 24. 03. 2021" data refresh to the Slovakia Covid
# daily stats sheet: updates cumulative AgTests (col F) and AgPosit (col G)
# counts for rows 334-383 (dates 2021-02-01 .. 2021-03-22), matching the
# upstream OpenData CSV refresh. Rows 354, 361 and 368 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(334, 6).Value = 196171
$ws.Cells.Item(334, 7).Value = 3479
$ws.Cells.Item(335, 6).Value = 131105
$ws.Cells.Item(335, 7).Value = 3007
$ws.Cells.Item(336, 6).Value = 101981
$ws.Cells.Item(336, 7).Value = 3389
$ws.Cells.Item(337, 6).Value = 104140
$ws.Cells.Item(337, 7).Value = 2961
$ws.Cells.Item(338, 6).Value = 227201
$ws.Cells.Item(338, 7).Value = 3187
$ws.Cells.Item(339, 6).Value = 659424
$ws.Cells.Item(339, 7).Value = 5493
$ws.Cells.Item(340, 6).Value = 384184
$ws.Cells.Item(340, 7).Value = 3299
$ws.Cells.Item(341, 6).Value = 291458
$ws.Cells.Item(341, 7).Value = 3665
$ws.Cells.Item(342, 6).Value = 179648
$ws.Cells.Item(342, 7).Value = 3072
$ws.Cells.Item(343, 6).Value = 132647
$ws.Cells.Item(343, 7).Value = 2968
$ws.Cells.Item(344, 6).Value = 135637
$ws.Cells.Item(345, 6).Value = 292028
$ws.Cells.Item(346, 6).Value = 671047
$ws.Cells.Item(346, 7).Value = 4787
$ws.Cells.Item(347, 6).Value = 342415
$ws.Cells.Item(348, 6).Value = 232280
$ws.Cells.Item(348, 7).Value = 3252
$ws.Cells.Item(349, 6).Value = 159098
$ws.Cells.Item(349, 7).Value = 2754
$ws.Cells.Item(350, 6).Value = 127068
$ws.Cells.Item(350, 7).Value = 2780
$ws.Cells.Item(351, 6).Value = 150624
$ws.Cells.Item(351, 7).Value = 2826
$ws.Cells.Item(352, 6).Value = 307175
$ws.Cells.Item(352, 7).Value = 3541
$ws.Cells.Item(353, 6).Value = 723222
$ws.Cells.Item(353, 7).Value = 5266
$ws.Cells.Item(355, 6).Value = 222168
$ws.Cells.Item(355, 7).Value = 3456
$ws.Cells.Item(356, 6).Value = 159851
$ws.Cells.Item(356, 7).Value = 2878
$ws.Cells.Item(357, 6).Value = 138575
$ws.Cells.Item(357, 7).Value = 3028
$ws.Cells.Item(358, 6).Value = 157292
$ws.Cells.Item(358, 7).Value = 2599
$ws.Cells.Item(359, 6).Value = 320284
$ws.Cells.Item(359, 7).Value = 3337
$ws.Cells.Item(360, 6).Value = 746046
$ws.Cells.Item(360, 7).Value = 5110
$ws.Cells.Item(362, 6).Value = 227808
$ws.Cells.Item(362, 7).Value = 3175
$ws.Cells.Item(363, 6).Value = 188331
$ws.Cells.Item(363, 7).Value = 2765
$ws.Cells.Item(364, 6).Value = 167336
$ws.Cells.Item(364, 7).Value = 2464
$ws.Cells.Item(365, 6).Value = 180421
$ws.Cells.Item(365, 7).Value = 2354
$ws.Cells.Item(366, 6).Value = 337794
$ws.Cells.Item(366, 7).Value = 2837
$ws.Cells.Item(367, 6).Value = 763050
$ws.Cells.Item(367, 7).Value = 3894
$ws.Cells.Item(369, 6).Value = 234314
$ws.Cells.Item(369, 7).Value = 2577
$ws.Cells.Item(370, 6).Value = 181583
$ws.Cells.Item(370, 7).Value = 2025
$ws.Cells.Item(371, 6).Value = 158152
$ws.Cells.Item(371, 7).Value = 1941
$ws.Cells.Item(372, 6).Value = 176460
$ws.Cells.Item(372, 7).Value = 1836
$ws.Cells.Item(373, 6).Value = 344154
$ws.Cells.Item(373, 7).Value = 2346
$ws.Cells.Item(374, 6).Value = 767541
$ws.Cells.Item(374, 7).Value = 3400
$ws.Cells.Item(375, 6).Value = 349231
$ws.Cells.Item(375, 7).Value = 1843
$ws.Cells.Item(376, 6).Value = 219547
$ws.Cells.Item(376, 7).Value = 2197
$ws.Cells.Item(377, 6).Value = 174957
$ws.Cells.Item(377, 7).Value = 1816
$ws.Cells.Item(378, 6).Value = 155039
$ws.Cells.Item(378, 7).Value = 1519
$ws.Cells.Item(379, 6).Value = 176267
$ws.Cells.Item(379, 7).Value = 1595
$ws.Cells.Item(380, 6).Value = 337369
$ws.Cells.Item(380, 7).Value = 1962
$ws.Cells.Item(381, 6).Value = 725886
$ws.Cells.Item(381, 7).Value = 2620
$ws.Cells.Item(382, 6).Value = 348212
$ws.Cells.Item(382, 7).Value = 1577
$ws.Cells.Item(383, 6).Value = 213502
$ws.Cells.Item(383, 7).Value = 1708
